# Add the "North Dakota (2024)" real-value column (E) and the formula
# parameter labels column (G) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen the new columns to match the authored layout (closest value the
# engine's column-width rounding can reproduce for the saved OOXML widths
# 20.33203125 / 26.6640625).
$ws.Columns.Item(5).ColumnWidth = 19.5
$ws.Columns.Item(7).ColumnWidth = 25.833333333333332

# Parameter labels in column G.
$ws.Range("G2").Value = "A=1500"
$ws.Range("G3").Value = "B=30/20"
$ws.Range("G4").Value = "C=8"
$ws.Range("G5").Value = "alpha = 0.95"
$ws.Range("G6").Value = "lambda = 0.5"

# Header for the new "real value" column.
$ws.Range("E1").Value = "North Dakota (2024)"

# RAINFALL block (row 2).
$ws.Range("E2").Value = 17.8

# Crop_price block real values (rows 5-11, CORN..HAY); COTTON (row 12) has no value.
$ws.Range("E5").Value = 4.21
$ws.Range("E6").Value = 9.44
$ws.Range("E7").Value = 5.62
$ws.Range("E8").Value = 3.4
$ws.Range("E9").Value = 5.35
$ws.Range("E11").Value = 89

# COTTON (row 12 for price) has no real value; leave E10 blank but present
# (mirrors the author's empty-but-touched cell) without minting a new style.
$ws.Range("E10").Borders.LineStyle = -4142

# Crop_yield block real values (rows 15-19, CORN..BARLEY).
$ws.Range("E15").Value = 149
$ws.Range("E16").Value = 37.5
$ws.Range("E17").Value = 50.5
$ws.Range("E18").Value = 98
$ws.Range("E19").Value = 74
$ws.Range("E21").Value = 1.77

# Match the author's final selection / active cell.
$ws.Range("G20").Select()
